$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "DataHazards Truth Table"

$newSheet.Range("A1").NumberFormat = "@"
$newSheet.Range("A1").Value = "if(A) C=10"
$newSheet.Range("A2").NumberFormat = "@"
$newSheet.Range("A2").Value = "else if(B) C=01"
$newSheet.Range("A3").NumberFormat = "@"
$newSheet.Range("A3").Value = "else C=00"

$newSheet.Range("A4:C4").NumberFormat = "@"
$newSheet.Range("A4").Value = "A"
$newSheet.Range("B4").Value = "B"
$newSheet.Range("C4").Value = "C"

$newSheet.Range("A5:C8").NumberFormat = "@"
$newSheet.Range("A5").Value = "0"
$newSheet.Range("B5").Value = "0"
$newSheet.Range("C5").Value = "00"

$newSheet.Range("A6").Value = "0"
$newSheet.Range("B6").Value = "1"
$newSheet.Range("C6").Value = "01"

$newSheet.Range("A7").Value = "1"
$newSheet.Range("B7").Value = "0"
$newSheet.Range("C7").Value = "10"

$newSheet.Range("A8").Value = "1"
$newSheet.Range("B8").Value = "1"
$newSheet.Range("C8").Value = "10"

$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

$newSheet.Range("A9").Select() | Out-Null
